$wb = $excel.ActiveWorkbook

# --- "Agile Test Plan" sheet updates ---
$ws = $wb.Worksheets.Item("Agile Test Plan")

# Row 21: new test case (Test #4 - set max price to 27000)
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = "Set maximum price to 27000"
$ws.Range("F21").Value = "Only Data entries with cost <= 27000 are displayed"
$ws.Range("G21").Value = "As expected"
$ws.Range("I21").Value = "Afzal"
$ws.Range("J21").Value = "Potential issue with labels leading to duplication, can maybe be resolved when map data is fully implemented. "

# Row 22: new test case (Test #5 - update price filter to 29000)
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = "Updating the price filter from 27000 to 29000"
$ws.Range("F22").Value = "update display to show Data entries with cost <= 29000"
$ws.Range("G22").Value = "As expected"
$ws.Range("I22").Value = "Afzal"

# Row 17: finish existing test case (maximum price -100)
$ws.Range("F17").Value = "No results found"
$ws.Range("G17").Value = "as expected"
$ws.Range("I17").Value = "Afzal"
$ws.Range("J17").Value = "Maybe revisit code to prevent negative numbers as oppopsed to searchhing with them and returning no values."
$ws.Rows.Item(17).RowHeight = 43.5

# --- "-Disclaimer-" sheet updates: mark additional checkboxes agreed/ticked ---
$wsD = $wb.Worksheets.Item("-Disclaimer-")
$wsD.Range("A19").Value = $true
$wsD.Range("A23").Value = $true
$wsD.Range("A24").Value = $true
